$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original row 13 (B13/C13 = "7926291 - Célia Regina Tomachuk dos Santos Catuogno",
# no A13 value) is removed entirely; everything below shifts up by one row.
$ws.Rows.Item(13).Delete()

# Row 10 (Objetivos:) — B/C body text replaced with the instructor line.
$ws.Range("B10").Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"
$ws.Range("C10").Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"

# Row 13 (now "Programa resumido:") — B/C replaced with "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (now "Programa:") — B/C replaced with the activation date string.
# Leading apostrophe forces literal text so Excel doesn't reinterpret the
# date-shaped string as a serial date value.
$ws.Range("B15").Formula = "'01/01/2012"
$ws.Range("C15").Formula = "'01/01/2012"

# Row 18 (now "Método:") — B/C replaced with the instructor line again.
$ws.Range("B18").Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"
$ws.Range("C18").Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"

# Row 19 (now "Critério:") — B/C replaced with the lecture-method text.
$ws.Range("B19").Value = "Aulas expositivas, exercícios e visitas didátias de campo."
$ws.Range("C19").Value = "Aulas expositivas, exercícios e visitas didátias de campo."

# Row 20 (now "Norma de recuperação:") — B/C replaced with the grading-criteria text.
$ws.Range("B20").Value = "Provas e relatórios."
$ws.Range("C20").Value = "Provas e relatórios."

# Row 21 (now "Bibliografia:") — B/C replaced with the single-exam-pass text.
$ws.Range("B21").Value = "Prova única com nota igual ou superior a 5,0 (cinco)."
$ws.Range("C21").Value = "Prova única com nota igual ou superior a 5,0 (cinco)."
